$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "i79x1a05b5"
$ws.Range("B12").Select()
